# "Generate Report for Archive"
# - Update localization status text from "Ready for handoff" to "In Translation"
#   on every sheet that references it (Overview!E2/F2, zh-cn!C2, de-de!C2).
# - Columns affected by the shorter status text re-flow to a narrower width
#   on each of those sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Update the status value everywhere it appears so the shared string is
# rewritten in place rather than leaving the old text orphaned.
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# Narrow the affected status columns to match the new content width.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
